$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 141, shifting existing rows 141:230 down to 142:231
$ws.Rows.Item(141).Insert()

# Fill in the new row 141 with values (same template as surrounding rows, but new data)
$ws.Range("A141").Value = 4
$ws.Range("B141").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C141").Value = "Los Lagos"
$ws.Range("D141").Value = 44596
$ws.Range("D141").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E141").Value = 10
$ws.Range("F141").Value = 100112037
$ws.Range("G141").Value = "Cebollín"
$ws.Range("H141").Value = "Sin especificar"
$ws.Range("I141").Value = "Primera"
$ws.Range("J141").Value = 180
$ws.Range("K141").Value = 5000
$ws.Range("L141").Value = 5500
$ws.Range("M141").Value = 5250
$ws.Range("N141").Value = "$/paquete 36 unidades"
$ws.Range("O141").Value = "Región Metropolitana"
$ws.Range("P141").Value = 146
$ws.Range("Q141").Value = 36
$ws.Range("R141").Value = "Hortaliza"
